$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.50"
$ws.Range("D4").Value = "'5.397"
$ws.Range("D5").Value = "'0.05759"
$ws.Range("D6").Value = "'3.435"
$ws.Range("D7").Value = "'6.312"
$ws.Range("D8").Value = "'0.8106"
$ws.Range("D9").Value = "'0.8928"
$ws.Range("D10").Value = "'0.1442"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("D11").Value = "'0.07336"
$ws.Range("D14").Value = "'0.09420"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001576"
$ws.Range("E15").Value = "14BitForexTokenBF"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").Value = "'0.04792"
$ws.Range("E16").Value = "15CoinExTokenCET"
$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D17").Value = "'0.0005843"
$ws.Range("E17").Value = "16OneONE"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "'0.006320"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D19").Value = "'0.004064"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").Value = "'0.0009948"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").Value = "'0.0001501"
$ws.Range("E21").Value = "20NitroExNTX"
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").Value = "'3.727"
$ws.Range("E22").Value = "21LEOLEO"
$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D23").Value = "'2.192"
$ws.Range("E23").Value = "22BTSETokenBTSE"
$ws.Range("B24").Value = "BitpandaEcosystemToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D24").Value = "'0.3274"
$ws.Range("E24").Value = "23BitpandaEcosystemTokenBEST"
$ws.Range("B25").Value = "ProBitToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D25").Value = "'0.1301"
$ws.Range("E25").Value = "24ProBitTokenPROB"
$ws.Range("B26").Value = "MCDex"
$ws.Range("C26").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D26").Value = "'4.164"
$ws.Range("E26").Value = "25MCDexMCB"
$ws.Range("D27").Value = "'0.0004652"
$ws.Range("D40").Value = "'0.03892"
$ws.Range("D41").Value = "'0.006779"
$ws.Range("D42").Value = "'0.1071"
$ws.Range("D43").Value = "'0.003202"
$ws.Range("E43").Value = "42CEJICEJIBestin24h"
$ws.Range("D45").Value = "'0.00005643"
$ws.Range("D47").Value = "'0.3802"
$ws.Range("D48").Value = "'0.1640"
$ws.Range("D50").Value = "'0.01011"
